$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "H358-2"
$ws.Range("A5").Value = "H2122-2"
$ws.Range("A1").Value = "Barcode"

$ws.Range("C2").Value = "D300_trt_Tecan_96_well_plates.xlsx"
$ws.Range("C3").Value = "D300_trt_Tecan_96_well_plates.xlsx"
$ws.Range("C4").Value = "D300_trt_Tecan_96_well_plates.xlsx"
$ws.Range("C5").Value = "D300_trt_Tecan_96_well_plates.xlsx"

$ws.Columns.Item(3).ColumnWidth = 31

$ws.Range("C5").Select()
